# Regenerate vocabulary sheet from updated Google Sheet export:
#  - widen used range to columns AL:AM (new trailing blank columns)
#  - add new rows up to 29
#  - bump "modified" timestamp (B20)
#  - insert a new "skos:broadMatch" header column before "iop:hasProperty" (row 22),
#    shifting the following headers one column to the right
#  - flesh out the two placeholder terms (rows 23-24) with real labels + provenance
#  - append new terms vocab:1002 .. vocab:1006 (rows 25-29)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 0. Widen the used range first (AL:AM, through row 29) -----------------
# A plain NumberFormat touch on the whole target block registers the blank
# cells (and new rows) without putting any real value in them, which is what
# the source diff shows for these columns: present, but empty.
$ws.Range("AL1:AM29").NumberFormat = "General"
# New rows 25-29 are (like every other row in this export) fully populated
# A:AM with either real text or an empty placeholder cell; pre-touch the
# whole block so the untouched columns still materialise as blank cells.
$ws.Range("A25:AM29").NumberFormat = "General"

# --- 1. dct:modified^^xsd:datetime bump -------------------------------------
$ws.Range("B20").Value = "2022-06-23T15:16:00+00:00"

# --- 2. Row 22 header row: insert skos:broadMatch before iop:hasProperty ---
# (values I22..W22 all shift right by one into J22..X22)
$ws.Range("I22").Value = 'skos:broadMatch(separator=",")'
$ws.Range("J22").Value = "iop:hasProperty"
$ws.Range("K22").Value = "iop:hasObjectOfInterest"
$ws.Range("L22").Value = "iop:hasMatrix"
$ws.Range("M22").Value = 'iop:hasContextObject(separator=",")'
$ws.Range("N22").Value = 'iop:hasConstraint(separator=",")'
$ws.Range("O22").Value = 'puv:statistic(separator=",")'
$ws.Range("P22").Value = 'puv:usesMethod(separator=",")'
$ws.Range("Q22").Value = 'sosa:madeBySensor(separator=",")'
$ws.Range("R22").Value = 'puv:uom(separator=",")'
$ws.Range("S22").Value = "owl:deprecated^^xsd:boolean"
$ws.Range("T22").Value = "skos:editorialNote@en"
$ws.Range("U22").Value = "dct:modified^^xsd:date"
$ws.Range("V22").Value = "dct:created^^xsd:date"
$ws.Range("W22").Value = 'dct:creator(separator=",")'
$ws.Range("X22").Value = 'dct:contributor(separator=",")'

# --- helper: write a cell as literal text even if it looks like a date -----
function Set-TextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

# --- 3. Row 23 (vocab:1000) — fill in real label + provenance --------------
$ws.Range("B23").Value = "performance of feature"
Set-TextValue "U23" "2022-06-23"
Set-TextValue "V23" "2022-06-23"
$ws.Range("W23").Value = "0000-0003-2195-3997"

# --- 4. Row 24 (vocab:1001) — fill in real label + broader + provenance ----
$ws.Range("B24").Value = "pre-surprise performance of key feature"
$ws.Range("F24").Value = "performance of feature"
Set-TextValue "U24" "2022-06-23"
Set-TextValue "V24" "2022-06-23"
$ws.Range("W24").Value = "0000-0003-2195-3997"

# --- 5. New rows 25-29 ------------------------------------------------------
$ws.Range("A25").Value = "vocab:1002"
$ws.Range("B25").Value = "surprise performance of surprise feature"
$ws.Range("F25").Value = "performance of feature"
Set-TextValue "U25" "2022-06-23"
Set-TextValue "V25" "2022-06-23"
$ws.Range("W25").Value = "0000-0003-2195-3997"

$ws.Range("A26").Value = "vocab:1003"
$ws.Range("B26").Value = "post-surprise performance of key feature"
$ws.Range("F26").Value = "performance of feature"
Set-TextValue "U26" "2022-06-23"
Set-TextValue "V26" "2022-06-23"
$ws.Range("W26").Value = "0000-0003-2195-3997"

$ws.Range("A27").Value = "vocab:1004"
$ws.Range("B27").Value = "post-surprise performance of surprise feature"
$ws.Range("F27").Value = "performance of feature"
Set-TextValue "U27" "2022-06-23"
Set-TextValue "V27" "2022-06-23"
$ws.Range("W27").Value = "0000-0003-2195-3997"

$ws.Range("A28").Value = "vocab:1005"
$ws.Range("B28").Value = "spectral power"
$ws.Range("D28").Value = "Spectral power of the signal measured from the data."
$ws.Range("E28").Value = "https://www.sciencedirect.com/topics/engineering/power-spectrum"
$ws.Range("O28").Value = "http://purl.bioontology.org/ontology/SNOMEDCT/255586005"
$ws.Range("P28").Value = "http://bioontology.org/ontologies/BiomedicalResourceOntology.owl#Fourier_Transform"
Set-TextValue "U28" "2022-06-23"
Set-TextValue "V28" "2022-06-23"
$ws.Range("W28").Value = "0000-0001-6361-2571"

$ws.Range("A29").Value = "vocab:1006"
$ws.Range("B29").Value = "inter-trial coherence"
$ws.Range("D29").Value = "Inter-trial coherence between various trials (epochs) across a measurement. "
$ws.Range("E29").Value = "https://pubmed.ncbi.nlm.nih.gov/24360131/"
$ws.Range("O29").Value = "http://purl.bioontology.org/ontology/SNOMEDCT/255586005"
Set-TextValue "U29" "2022-06-23"
Set-TextValue "V29" "2022-06-23"
$ws.Range("W29").Value = "0000-0001-6361-2571"

Write-Output "edit complete"
